$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3.182878228561681
$ws.Range("C2").Value = 1.65323645889881
$ws.Range("D2").Value = 3.082599426703578
$ws.Range("E2").Value = 0.4998867070740569
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 8.418600821238126
$ws.Range("B3").Value = 1.505614041169197
$ws.Range("C3").Value = 1.65323645889881
$ws.Range("D3").Value = 0.1529057820181812
$ws.Range("E3").Value = 0.4998867070740569
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 3.811642989160245
$ws.Range("B4").Value = 3.182878228561681
$ws.Range("C4").Value = 1.65323645889881
$ws.Range("D4").Value = 3.082599426703578
$ws.Range("E4").Value = 0.4998867070740569
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 8.418600821238126
$ws.Range("B5").Value = 0.7287194209349384
$ws.Range("C5").Value = 0.05231270169004087
$ws.Range("D5").Value = 0.1529057820181812
$ws.Range("E5").Value = 0.4998867070740569
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.433824611717217
$ws.Range("B6").Value = 3.182878228561681
$ws.Range("C6").Value = 1.65323645889881
$ws.Range("D6").Value = 0.7127328510149897
$ws.Range("E6").Value = 0.4998867070740569
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 6.048734245549538
$ws.Range("B7").Value = 3.182878228561681
$ws.Range("C7").Value = 1.65323645889881
$ws.Range("D7").Value = 0.7127328510149897
$ws.Range("E7").Value = 0.4998867070740569
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 6.048734245549538
$ws.Range("B8").Value = 3.182878228561681
$ws.Range("C8").Value = 1.65323645889881
$ws.Range("D8").Value = 0.1529057820181812
$ws.Range("E8").Value = 0.4998867070740569
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 5.488907176552729
$ws.Range("B9").Value = 3.182878228561681
$ws.Range("C9").Value = 1.65323645889881
$ws.Range("D9").Value = 3.082599426703578
$ws.Range("E9").Value = 0.4998867070740569
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 8.418600821238126
$ws.Range("B10").Value = 3.182878228561681
$ws.Range("C10").Value = 1.65323645889881
$ws.Range("D10").Value = 0.7127328510149897
$ws.Range("E10").Value = 0.4998867070740569
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 6.048734245549538
$ws.Range("B11").Value = 1.505614041169197
$ws.Range("C11").Value = 1.65323645889881
$ws.Range("D11").Value = 0.7127328510149897
$ws.Range("E11").Value = 0.4998867070740569
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 4.371470058157054
$ws.Range("B12").Value = 3.182878228561681
$ws.Range("C12").Value = 1.65323645889881
$ws.Range("D12").Value = 0.1529057820181812
$ws.Range("E12").Value = 0.4998867070740569
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 5.488907176552729
$ws.Range("B13").Value = 1.505614041169197
$ws.Range("C13").Value = 1.65323645889881
$ws.Range("D13").Value = 0.7127328510149897
$ws.Range("E13").Value = 0.4998867070740569
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 4.371470058157054
$ws.Range("B14").Value = 3.182878228561681
$ws.Range("C14").Value = 1.65323645889881
$ws.Range("D14").Value = 3.082599426703578
$ws.Range("E14").Value = 0.4998867070740569
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 8.418600821238126
$ws.Range("B15").Value = 0.3464964993005633
$ws.Range("C15").Value = 0.05231270169004087
$ws.Range("D15").Value = 3.082599426703578
$ws.Range("E15").Value = 0.4998867070740569
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 3.981295334768239
$ws.Range("B16").Value = 0.06328177979961902
$ws.Range("C16").Value = 0.3375848360084654
$ws.Range("D16").Value = 0.7127328510149897
$ws.Range("E16").Value = 0.4998867070740569
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 1.613486173897131
$ws.Range("B17").Value = 0.1554434735375247
$ws.Range("C17").Value = 0.3375848360084654
$ws.Range("D17").Value = 0.1529057820181812
$ws.Range("E17").Value = 0.4998867070740569
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 1.145820798638228
$ws.Range("B18").Value = 0.3464964993005633
$ws.Range("C18").Value = 0.3375848360084654
$ws.Range("D18").Value = 0.1529057820181812
$ws.Range("E18").Value = 0.4998867070740569
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 1.336873824401267
$ws.Range("B19").Value = 0.02258322285507441
$ws.Range("C19").Value = 0.05231270169004087
$ws.Range("D19").Value = 0.7127328510149897
$ws.Range("E19").Value = 0.4998867070740569
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 1.287515482634162
$ws.Range("B20").Value = 3.182878228561681
$ws.Range("C20").Value = 1.65323645889881
$ws.Range("D20").Value = 0.7127328510149897
$ws.Range("E20").Value = 0.4998867070740569
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 6.048734245549538
$ws.Range("B21").Value = 3.182878228561681
$ws.Range("C21").Value = 1.65323645889881
$ws.Range("D21").Value = 0.7127328510149897
$ws.Range("E21").Value = 0.4998867070740569
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 6.048734245549538
$ws.Range("B22").Value = 0.7287194209349384
$ws.Range("C22").Value = 1.65323645889881
$ws.Range("D22").Value = 0.1529057820181812
$ws.Range("E22").Value = 0.4998867070740569
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 3.034748368925986
$ws.Range("B23").Value = 1.505614041169197
$ws.Range("C23").Value = 1.65323645889881
$ws.Range("D23").Value = 3.082599426703578
$ws.Range("E23").Value = 0.4998867070740569
$ws.Range("F23").Value = 1
$ws.Range("G23").Value = 6.741336633845642
$ws.Range("B24").Value = 3.182878228561681
$ws.Range("C24").Value = 1.65323645889881
$ws.Range("D24").Value = 0.7127328510149897
$ws.Range("E24").Value = 0.4998867070740569
$ws.Range("F24").Value = 1
$ws.Range("G24").Value = 6.048734245549538
$ws.Range("B25").Value = 1.505614041169197
$ws.Range("C25").Value = 1.65323645889881
$ws.Range("D25").Value = 3.082599426703578
$ws.Range("E25").Value = 0.4998867070740569
$ws.Range("F25").Value = 1
$ws.Range("G25").Value = 6.741336633845642
$ws.Range("B26").Value = 1.505614041169197
$ws.Range("C26").Value = 1.65323645889881
$ws.Range("D26").Value = 0.1529057820181812
$ws.Range("E26").Value = 0.4998867070740569
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = 3.811642989160245
$ws.Range("B27").Value = 1.505614041169197
$ws.Range("C27").Value = 1.65323645889881
$ws.Range("D27").Value = 0.1529057820181812
$ws.Range("E27").Value = 0.4998867070740569
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 3.811642989160245
$ws.Range("B28").Value = 3.182878228561681
$ws.Range("C28").Value = 1.65323645889881
$ws.Range("D28").Value = 16.98373111632243
$ws.Range("E28").Value = 0.4998867070740569
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 22.31973251085698
$ws.Range("B29").Value = 3.182878228561681
$ws.Range("C29").Value = 1.65323645889881
$ws.Range("D29").Value = 3.082599426703578
$ws.Range("E29").Value = 0.4998867070740569
$ws.Range("F29").Value = 1
$ws.Range("G29").Value = 8.418600821238126
$ws.Range("B30").Value = 3.182878228561681
$ws.Range("C30").Value = 1.65323645889881
$ws.Range("D30").Value = 0.7127328510149897
$ws.Range("E30").Value = 0.4998867070740569
$ws.Range("F30").Value = 1
$ws.Range("G30").Value = 6.048734245549538
$ws.Range("B31").Value = 3.182878228561681
$ws.Range("C31").Value = 1.65323645889881
$ws.Range("D31").Value = 3.082599426703578
$ws.Range("E31").Value = 6.48142807727062
$ws.Range("F31").Value = 1
$ws.Range("G31").Value = 14.40014219143469
$ws.Range("B32").Value = 1.505614041169197
$ws.Range("C32").Value = 1.65323645889881
$ws.Range("D32").Value = 0.1529057820181812
$ws.Range("E32").Value = 0.4998867070740569
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 3.811642989160245
$ws.Range("B33").Value = 1.505614041169197
$ws.Range("C33").Value = 1.65323645889881
$ws.Range("D33").Value = 3.082599426703578
$ws.Range("E33").Value = 0.4998867070740569
$ws.Range("F33").Value = 1
$ws.Range("G33").Value = 6.741336633845642
$ws.Range("B34").Value = 3.182878228561681
$ws.Range("C34").Value = 1.65323645889881
$ws.Range("D34").Value = 0.1529057820181812
$ws.Range("E34").Value = 0.4998867070740569
$ws.Range("F34").Value = 1
$ws.Range("G34").Value = 5.488907176552729
$ws.Range("B35").Value = 3.182878228561681
$ws.Range("C35").Value = 1.65323645889881
$ws.Range("D35").Value = 0.1529057820181812
$ws.Range("E35").Value = 0.4998867070740569
$ws.Range("F35").Value = 1
$ws.Range("G35").Value = 5.488907176552729
$ws.Range("B36").Value = 3.182878228561681
$ws.Range("C36").Value = 1.65323645889881
$ws.Range("D36").Value = 16.98373111632243
$ws.Range("E36").Value = 6.48142807727062
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 28.30127388105354
$ws.Range("B37").Value = 3.182878228561681
$ws.Range("C37").Value = 0.3375848360084654
$ws.Range("D37").Value = 0.1529057820181812
$ws.Range("E37").Value = 0.4998867070740569
$ws.Range("F37").Value = 1
$ws.Range("G37").Value = 4.173255553662385
$ws.Range("B38").Value = 0.7287194209349384
$ws.Range("C38").Value = 1.65323645889881
$ws.Range("D38").Value = 0.7127328510149897
$ws.Range("E38").Value = 0.4998867070740569
$ws.Range("F38").Value = 1
$ws.Range("G38").Value = 3.594575437922795
$ws.Range("B39").Value = 3.182878228561681
$ws.Range("C39").Value = 1.65323645889881
$ws.Range("D39").Value = 0.7127328510149897
$ws.Range("E39").Value = 0.4998867070740569
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = 6.048734245549538
$ws.Range("B40").Value = 3.182878228561681
$ws.Range("C40").Value = 1.65323645889881
$ws.Range("D40").Value = 16.98373111632243
$ws.Range("E40").Value = 0.4998867070740569
$ws.Range("F40").Value = 1
$ws.Range("G40").Value = 22.31973251085698
$ws.Range("B41").Value = 0.3464964993005633
$ws.Range("C41").Value = 0.3375848360084654
$ws.Range("D41").Value = 0.1529057820181812
$ws.Range("E41").Value = 0.4998867070740569
$ws.Range("F41").Value = 0
$ws.Range("G41").Value = 1.336873824401267
$ws.Range("B42").Value = 0.3464964993005633
$ws.Range("C42").Value = 1.65323645889881
$ws.Range("D42").Value = 0.1529057820181812
$ws.Range("E42").Value = 0.4998867070740569
$ws.Range("F42").Value = 1
$ws.Range("G42").Value = 2.652525447291612
$ws.Range("B43").Value = 3.182878228561681
$ws.Range("C43").Value = 1.65323645889881
$ws.Range("D43").Value = 0.7127328510149897
$ws.Range("E43").Value = 0.4998867070740569
$ws.Range("F43").Value = 1
$ws.Range("G43").Value = 6.048734245549538
$ws.Range("B44").Value = 3.182878228561681
$ws.Range("C44").Value = 1.65323645889881
$ws.Range("D44").Value = 0.7127328510149897
$ws.Range("E44").Value = 0.4998867070740569
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 6.048734245549538
$ws.Range("B45").Value = 3.182878228561681
$ws.Range("C45").Value = 1.65323645889881
$ws.Range("D45").Value = 0.7127328510149897
$ws.Range("E45").Value = 0.4998867070740569
$ws.Range("F45").Value = 1
$ws.Range("G45").Value = 6.048734245549538
$ws.Range("B46").Value = 3.182878228561681
$ws.Range("C46").Value = 1.65323645889881
$ws.Range("D46").Value = 0.7127328510149897
$ws.Range("E46").Value = 0.4998867070740569
$ws.Range("F46").Value = 1
$ws.Range("G46").Value = 6.048734245549538
$ws.Range("B47").Value = 3.182878228561681
$ws.Range("C47").Value = 1.65323645889881
$ws.Range("D47").Value = 0.7127328510149897
$ws.Range("E47").Value = 0.4998867070740569
$ws.Range("F47").Value = 1
$ws.Range("G47").Value = 6.048734245549538
$ws.Range("B48").Value = 3.182878228561681
$ws.Range("C48").Value = 1.65323645889881
$ws.Range("D48").Value = 16.98373111632243
$ws.Range("E48").Value = 0.4998867070740569
$ws.Range("F48").Value = 1
$ws.Range("G48").Value = 22.31973251085698
$ws.Range("B49").Value = 1.505614041169197
$ws.Range("C49").Value = 1.65323645889881
$ws.Range("D49").Value = 3.082599426703578
$ws.Range("E49").Value = 0.4998867070740569
$ws.Range("F49").Value = 1
$ws.Range("G49").Value = 6.741336633845642
$ws.Range("B50").Value = 0.1554434735375247
$ws.Range("C50").Value = 0.3375848360084654
$ws.Range("D50").Value = 0.1529057820181812
$ws.Range("E50").Value = 0.4998867070740569
$ws.Range("F50").Value = 0
$ws.Range("G50").Value = 1.145820798638228
